# Renumber the village scene names so the first scene ("villageScene")
# is no longer an empty/zero-suffixed name - the protocol-body length
# must be > 0, so every scene name below gets shifted up by one and a
# new trailing "villageScene6" entry is introduced for the last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "villageScene1"
$ws.Range("C11").Value = "villageScene1"

$ws.Range("B12").Value = "villageScene2"
$ws.Range("C12").Value = "villageScene2"

$ws.Range("B13").Value = "villageScene3"
$ws.Range("C13").Value = "villageScene3"

$ws.Range("B14").Value = "villageScene4"
$ws.Range("C14").Value = "villageScene4"

$ws.Range("B15").Value = "villageScene5"
$ws.Range("C15").Value = "villageScene5"

$ws.Range("B16").Value = "villageScene6"
$ws.Range("C16").Value = "villageScene6"
